$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("./img/player_image_20240713-234946.png", 24),
    @("./img/player_image_20240714-000050.png", 26),
    @("./img/player_image_20240714-000341.png", 14),
    @("./img/player_image_20240714-001913.png", 9),
    @("./img/player_image_20240714-002019.png", 10),
    @("./img/player_image_20240714-002138.png", 17),
    @("./img/player_image_20240714-002309.png", 7),
    @("./img/player_image_20240714-002851.png", 23),
    @("./img/player_image_20240714-005014.png", 4),
    @("./img/player_image_20240714-005124.png", 3),
    @("./img/player_image_20240714-005320.png", 12),
    @("./img/player_image_20240714-005449.png", 17),
    @("./img/player_image_20240714-005816.png", 9),
    @("./img/player_image_20240714-010859.png", 11),
    @("./img/player_image_20240714-011315.png", 12),
    @("./img/player_image_20240714-011814.png", 6)
)

$startRow = 10
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
